# Apply updates to column F (dSF) values for rows 2,3,5,6,7
# per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 3

$wb.Save()
